$d = $word.ActiveDocument
$d.Content.Find.Execute("757÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "302÷6=", 2)
$d.Content.Find.Execute("758÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "463÷2=", 2)
$d.Content.Find.Execute("221÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "977÷3=", 2)
$d.Content.Find.Execute("608÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "362÷5=", 2)
$d.Content.Find.Execute("475÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "389÷7=", 2)
$d.Content.Find.Execute("359÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "999÷8=", 2)
$d.Content.Find.Execute("831÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "203÷2=", 2)
$d.Content.Find.Execute("775÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "114÷5=", 2)
$d.Content.Find.Execute("176÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "801÷8=", 2)
$d.Content.Find.Execute("599÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "332÷6=", 2)
$d.Content.Find.Execute("290÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "814÷9=", 2)
$d.Content.Find.Execute("586÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "795÷6=", 2)
$d.Content.Find.Execute("308÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "522÷9=", 2)
$d.Content.Find.Execute("873÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "227÷6=", 2)
$d.Content.Find.Execute("847÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "316÷5=", 2)
$d.Content.Find.Execute("968÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "842÷8=", 2)
$d.Content.Find.Execute("642÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "144÷7=", 2)
$d.Content.Find.Execute("191÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "592÷3=", 2)
$d.Content.Find.Execute("407÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "369÷3=", 2)
$d.Content.Find.Execute("183÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "226÷4=", 2)
$d.Content.Find.Execute("171÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "492÷8=", 2)
$d.Content.Find.Execute("154÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "120÷4=", 2)
$d.Content.Find.Execute("827÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "819÷7=", 2)
$d.Content.Find.Execute("873÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "775÷6=", 2)
$d.Content.Find.Execute("227÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "137÷4=", 2)
